$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - copy formatting from the existing H1 header cell
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data values for columns I and J
$iValues = @(3, 8, 2, 9, 9, 5, 5, 8, 6, 3, 6, 2)
$jValues = @(5, 9, 5, 9, 9, 9, 6, 8, 6, 3, 7, 2)

for ($r = 0; $r -lt $iValues.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
